$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Add two new rows to the Financial section (Bain Consulting, Accenture)
# right after "Discover" and before "McKinsey". Inserting rows 22:23 pushes
# McKinsey/Square/Plaid down by two rows and auto-extends/shifts the
# surrounding merged cells (A21:A25 -> A21:A27, A26:A30 -> A28:A32, A31:A33 -> A33:A35).
$ws.Rows("22:23").Insert()
$ws.Range("B22").Value2 = "Bain Consulting"
$ws.Range("B23").Value2 = "Accenture"

# --- Step 2: Insert a new "Positions" column before the old column D
# ("Ways to Connect"), shifting Ways to Connect / Application Method/Content /
# Timeline one column to the right (D->E, E->F, F->G).
$ws.Columns("D").Insert()
$ws.Columns("D").ColumnWidth = 41.7109375

# Header
$ws.Range("D1").Value2 = "Positions"

# Big Name section positions
$ws.Range("D6").Value2 = "Data Scientist "
$ws.Range("D7").Value2 = "Data Scientist, PM intern"
$ws.Range("E7").Value2 = "UCD alumni recruiters"

# Startups section positions
$ws.Range("D10").Value2 = "no intern application, email/linkedin directly"
$ws.Range("D11").Value2 = "Data Scientist Intern (3: algorithms analytics inference)"
$ws.Range("D12").Value2 = "none?"
$ws.Range("D13").Value2 = "none?"
$ws.Range("D14").Value2 = "none?"
$ws.Range("D15").Value2 = "logistics analyst"

# Leave the selection where the author last left it (the newly typed
# "Accenture" cell in the Financial section).
$ws.Range("B23").Select() | Out-Null
